$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Wie besprochen bearbeiten wir aktuell den Auftrag A0000001. Hierbei müssen wichtige Tätigkeiten erledigt werden. ", 1),
    @("Bezogen auf Auftrag A0000002 haben wir ein Problem. ", 1),
    @("Ich melde mich wegen Auftrag A0000003. Alles läuft super und sollte bald fertiggestellt sein", 1),
    @("Es geht um A0000004. Wir schaffen das so nicht und benötigen bessere Unterstützung. ", 1),
    @("Ich melde mich bezüglich A0000005 hat alles geklappt.", 1),
    @("Ich melde mich wegen Auftrag A0000006.", 1),
    @("Ich finde den Prozess im Auftrag A0000007 gut.", 0),
    @("Wegen Auftrag A0000008. Alles passt.", 1),
    @("Es geht um Auftrag A0000009.", 1),
    @("Ich schreibe wegen Auftrag A0000010.", 1),
    @("Bezüglich Auftrag A0000011.", 1),
    @("Nochmal wegen Auftrag A0000012.", 1),
    @("Ich melde mich wegen dem Auftrag A0000013.", 1),
    @("Ich melde mich wegen des Auftrages A0000014.", 1),
    @("Wegen dem Auftarg A0000015.", 1),
    @("Wir haben bereits im Auftrag A0000016 ähnliche Tätigkeiten erledigt.", 0),
    @(" Mit A0000017 stimmt etwas nicht.", 0),
    @("Bei A0000018 hat es doch auch geklappt.", 0),
    @("Ich hoffe gestern hat mit dem Auftrag A0000019 alles geklappt.", 0),
    @("Wir können die rechtzeitige fertigstellung nicht mehr garantieren da wir durch den Auftrag A0000020 aufgehalten werden.", 0),
    @("Hat bei Auftrag A0000021 alles geklappt?", 0)
)

$row = 1
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

$ws.Range("A22").Select() | Out-Null
